$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Yahoo")
$ws.Activate()
